$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 119.2
$ws.Range("I5").Value = 124
$ws.Range("K5").Value = 124
$ws.Range("M5").Value = -9
$ws.Range("H12").Value = 342
$ws.Range("I12").Value = 337.16666
$ws.Range("K12").Value = 337.16666
$ws.Range("M12").Value = -167.16666
$ws.Range("H15").Value = 949.45557
$ws.Range("I15").Value = 949.45557
$ws.Range("K15").Value = 2848.36671
$ws.Range("M15").Value = -2679.36671
$ws.Range("H34").Value = 3811.1428
$ws.Range("I34").Value = 3811.1428
$ws.Range("K34").Value = 3811.1428
$ws.Range("M34").Value = -3608.1428
$ws.Range("H36").Value = 3811.1428
$ws.Range("I36").Value = 3811.1428
$ws.Range("K36").Value = 3811.1428
$ws.Range("M36").Value = -3096.1428
$ws.Range("H86").Value = 10583733
$ws.Range("I86").Value = 3396.2
$ws.Range("K86").Value = 3396.2
$ws.Range("M86").Value = -2273.2
$ws.Range("H89").Value = 10583733
$ws.Range("I89").Value = 3396.2
$ws.Range("K89").Value = 16981
$ws.Range("M89").Value = -11365
$ws.Range("H112").Value = 1768.7
$ws.Range("J112").Value = 1768.7
$ws.Range("L112").Value = 5306.1
$ws.Range("N112").Value = -7522.1
$ws.Range("H125").Value = 1250
$ws.Range("J125").Value = 1250
$ws.Range("L125").Value = 11250
$ws.Range("N125").Value = -16170
$ws.Range("H127").Value = 2515.8333
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 2515.8333
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 7547.499899999999
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -17467.4999
$ws.Range("H138").Value = 3673.1516
$ws.Range("J138").Value = 3921.4736
$ws.Range("L138").Value = 11764.4208
$ws.Range("N138").Value = -22044.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1583.6923
$ws.Range("I2").Value = 1283.9
$ws.Range("J2").Value = 2583
$ws.Range("K2").Value = 1283.9
$ws.Range("L2").Value = 2583
$ws.Range("M2").Value = -1170.9
$ws.Range("N2").Value = -2809
$ws.Range("H74").Value = 4688.95
$ws.Range("I74").Value = 3437.5386
$ws.Range("K74").Value = 3437.5386
$ws.Range("M74").Value = -2563.5386
$ws.Range("H77").Value = 4688.95
$ws.Range("I77").Value = 3437.5386
$ws.Range("K77").Value = 17187.693
$ws.Range("M77").Value = -12819.693
$ws.Range("H116").Value = 1583.6923
$ws.Range("I116").Value = 1283.9
$ws.Range("J116").Value = 2583
$ws.Range("K116").Value = 1283.9
$ws.Range("L116").Value = 2583
$ws.Range("M116").Value = 1010.1
$ws.Range("N116").Value = -7171
$ws.Range("H122").Value = 3471.879
$ws.Range("I122").Value = 2068.842
$ws.Range("K122").Value = 6206.526
$ws.Range("M122").Value = -3756.526

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 13000.2
$ws.Range("I97").Value = 13000.2
$ws.Range("K97").Value = 13000.2
$ws.Range("M97").Value = -12009.2
$ws.Range("H107").Value = 45638852
$ws.Range("I107").Value = 252492.88
$ws.Range("J107").Value = 166669140
$ws.Range("K107").Value = 252492.88
$ws.Range("L107").Value = 166669140
$ws.Range("M107").Value = -250572.88
$ws.Range("N107").Value = -166672980
$ws.Range("H134").Value = 2633.1396
$ws.Range("I134").Value = 2142.1082
$ws.Range("K134").Value = 6426.3246
$ws.Range("M134").Value = -3891.3246
$ws.Range("H139").Value = 64897.5
$ws.Range("J139").Value = 79796
$ws.Range("L139").Value = 79796
$ws.Range("N139").Value = -90076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4032.4688
$ws.Range("I31").Value = 2868.5833
$ws.Range("J31").Value = 4730.8
$ws.Range("K31").Value = 2868.5833
$ws.Range("L31").Value = 4730.8
$ws.Range("M31").Value = -2573.5833
$ws.Range("N31").Value = -5320.8
$ws.Range("H34").Value = 4032.4688
$ws.Range("I34").Value = 2868.5833
$ws.Range("J34").Value = 4730.8
$ws.Range("K34").Value = 2868.5833
$ws.Range("L34").Value = 4730.8
$ws.Range("M34").Value = -2666.5833
$ws.Range("N34").Value = -5134.8
$ws.Range("H132").Value = 305352.1
$ws.Range("I132").Value = 1927.4783
$ws.Range("K132").Value = 5782.4349
$ws.Range("M132").Value = -3252.4349

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 10659.444
$ws.Range("I81").Value = 956.5
$ws.Range("J81").Value = 13431.714
$ws.Range("K81").Value = 2869.5
$ws.Range("L81").Value = 40295.142
$ws.Range("M81").Value = -1746.5
$ws.Range("N81").Value = -42541.142
$ws.Range("H84").Value = 10659.444
$ws.Range("I84").Value = 956.5
$ws.Range("J84").Value = 13431.714
$ws.Range("K84").Value = 8608.5
$ws.Range("L84").Value = 120885.426
$ws.Range("M84").Value = -2992.5
$ws.Range("N84").Value = -132117.426
$ws.Range("H132").Value = 2716.12
$ws.Range("I132").Value = 1524.7778
$ws.Range("J132").Value = 2878.5757
$ws.Range("K132").Value = 13723.0002
$ws.Range("L132").Value = 25907.1813
$ws.Range("M132").Value = -11193.0002
$ws.Range("N132").Value = -30967.1813
$ws.Range("H133").Value = 6668.125
$ws.Range("I133").Value = 4669
$ws.Range("K133").Value = 14007
$ws.Range("M133").Value = -8947
$ws.Range("H136").Value = 5747
$ws.Range("I136").Value = 5747
$ws.Range("K136").Value = 17241
$ws.Range("M136").Value = -12141
$ws.Range("H141").Value = 7569.857
$ws.Range("J141").Value = 2500
$ws.Range("L141").Value = 7500
$ws.Range("N141").Value = -17860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14274.556
$ws.Range("I70").Value = 11397
$ws.Range("K70").Value = 11397
$ws.Range("M70").Value = -11127
$ws.Range("H73").Value = 14274.556
$ws.Range("I73").Value = 11397
$ws.Range("K73").Value = 11397
$ws.Range("M73").Value = -10461
$ws.Range("H126").Value = 8995.036
$ws.Range("I126").Value = 15727.444
$ws.Range("K126").Value = 47182.33199999999
$ws.Range("M126").Value = -44712.33199999999
$ws.Range("H134").Value = 49994
$ws.Range("J134").Value = 49994
$ws.Range("L134").Value = 149982
$ws.Range("N134").Value = -155052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3373.3
$ws.Range("I68").Value = 3089.2307
$ws.Range("J68").Value = 3900.8572
$ws.Range("K68").Value = 3089.2307
$ws.Range("L68").Value = 3900.8572
$ws.Range("M68").Value = -2340.2307
$ws.Range("N68").Value = -5398.8572
$ws.Range("H71").Value = 3373.3
$ws.Range("I71").Value = 3089.2307
$ws.Range("J71").Value = 3900.8572
$ws.Range("K71").Value = 15446.1535
$ws.Range("L71").Value = 19504.286
$ws.Range("M71").Value = -11702.1535
$ws.Range("N71").Value = -26992.286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2652946.8
$ws.Range("J62").Value = 8500
$ws.Range("L62").Value = 8500
$ws.Range("N62").Value = -9748
$ws.Range("H65").Value = 2652946.8
$ws.Range("J65").Value = 8500
$ws.Range("L65").Value = 42500
$ws.Range("N65").Value = -48740
$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180
$ws.Range("H136").Value = 6815.5654
$ws.Range("I136").Value = 7697.148
$ws.Range("J136").Value = 5562.7896
$ws.Range("K136").Value = 23091.444
$ws.Range("L136").Value = 16688.3688
$ws.Range("M136").Value = -20541.444
$ws.Range("N136").Value = -21788.3688
